$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Location Name / Location ID values between row pairs
# (3,4), (7,8), (9,10) so Canby/Molalla rows are exchanged.
$rowPairs = @(@(3,4), @(7,8), @(9,10))

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $nameA = $ws.Cells.Item($r1, 1).Value2
    $idA   = $ws.Cells.Item($r1, 2).Value2
    $nameB = $ws.Cells.Item($r2, 1).Value2
    $idB   = $ws.Cells.Item($r2, 2).Value2

    $ws.Cells.Item($r1, 1).Value2 = $nameB
    $ws.Cells.Item($r1, 2).Value2 = $idB
    $ws.Cells.Item($r2, 1).Value2 = $nameA
    $ws.Cells.Item($r2, 2).Value2 = $idA
}
